$d = $word.ActiveDocument

# 1. "analyzed" -> "investigated" in the survey experiment bullet.
$d.Content.Find.Execute("Designed, fielded, and analyzed a survey experiment", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Designed, fielded, and investigated a survey experiment", 2)

# 2. Update "Spring 2021 - Spring 2023" (Computer Science Lab Assistant) to "Spring 2022 - Spring 2023".
$d.Content.Find.Execute("Spring 2021 " + [char]8211 + " Spring 2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Spring 2022 " + [char]8211 + " Spring 2023", 2)
